$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("Q4").Value = 2.08
$ws.Range("R4").Value = 1.73

# Row 5
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.9
$ws.Range("R5").Value = 1.95

# Row 8
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65

# Row 9
$ws.Range("G9").Value = 2.35
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 3.1
$ws.Range("L9").Value = 3.75
$ws.Range("Y9").Value = 9.5
$ws.Range("Z9").Value = 21
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 15
$ws.Range("AJ9").Value = 12
$ws.Range("AK9").Value = 34
$ws.Range("AL9").Value = 26
$ws.Range("AM9").Value = 41
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 13
$ws.Range("AQ9").Value = 41
$ws.Range("AW9").Value = 5

# Row 12
$ws.Range("O12").Value = 1.17
$ws.Range("P12").Value = 5

# Row 16
$ws.Range("G16").Value = 1.5
$ws.Range("H16").Value = 3.6
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 1.98
$ws.Range("K16").Value = 2.15
$ws.Range("L16").Value = 6.8
$ws.Range("M16").Value = 1.03
$ws.Range("N16").Value = 6.8
$ws.Range("O16").Value = 1.34
$ws.Range("P16").Value = 2.75
$ws.Range("Q16").Value = 1.98
$ws.Range("R16").Value = 1.65
$ws.Range("T16").Value = 2.55
$ws.Range("U16").Value = 2.07
$ws.Range("V16").Value = 1.6
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 6.2
$ws.Range("Y16").Value = 8.25
$ws.Range("Z16").Value = 10
$ws.Range("AA16").Value = 13
$ws.Range("AB16").Value = 32
$ws.Range("AC16").Value = 8
$ws.Range("AD16").Value = 7.3
$ws.Range("AE16").Value = 21
$ws.Range("AF16").Value = 120
$ws.Range("AH16").Value = 14.5
$ws.Range("AI16").Value = 45
$ws.Range("AJ16").Value = 23
$ws.Range("AK16").Value = 200
$ws.Range("AL16").Value = 100
$ws.Range("AM16").Value = 100
$ws.Range("AN16").Value = 3.15
$ws.Range("AO16").Value = 6.8
$ws.Range("AP16").Value = 17.5
$ws.Range("AQ16").Value = 21
$ws.Range("AS16").Value = 200
$ws.Range("AT16").Value = 2.5
$ws.Range("AU16").Value = 8.25
$ws.Range("AV16").Value = 90
$ws.Range("AW16").Value = 8.25
$ws.Range("AX16").Value = 45
$ws.Range("AY16").Value = 50
$ws.Range("AZ16").Value = 350
$ws.Range("BA16").Value = 400
